$d = $word.ActiveDocument

# The edited element is the floating text box ("Cuadro de texto 1",
# docPr id 550983624) that sits right after the "...de acuerdo con el
# cronograma siguiente:" sentence in CLAUSULA TERCERA. It currently has
# a white fill + thin black outline (a visible "marco"/frame around the
# schedule placeholder) and the commit removes that frame by switching
# the shape to no fill / no line, adding an explicit wps:style block,
# and refreshing the effectExtent / editId that Word regenerates
# whenever the shape's geometry-affecting properties change.

# Confirm the target shape is present before editing.
$shapeFound = $false
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    if ($d.Shapes.Item($i).Id -eq 550983624) { $shapeFound = $true }
}
if (-not $shapeFound) {
    throw "Expected shape id 550983624 not present"
}

# Locate the paragraph that hosts the drawing: it is the very next
# paragraph after the "...cronograma siguiente:" sentence.
$anchorText = "cronograma siguiente:"
$matchRange = $d.Content
$found = $matchRange.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the anchor sentence before the schedule text box"
}
$hostStart = $matchRange.End + 1

$probe = $d.Range($hostStart, $hostStart + 1)
$hostPara = $probe.Paragraphs.Item(1)
$hostRange = $d.Range($hostPara.Range.Start, $hostPara.Range.End)

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4F8C259A" w14:textId="6EF22FB4" w:rsidR="00E34B0C" w:rsidRPr="000D1991" w:rsidRDefault="00CB7522" w:rsidP="00E34B0C"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr><w:rFonts w:ascii="Arial Rounded MT Bold" w:hAnsi="Arial Rounded MT Bold"/><w:lang w:val="es-PE"/></w:rPr></w:pPr><w:r w:rsidRPr="000D1991"><w:rPr><w:rFonts w:ascii="Arial Rounded MT Bold" w:hAnsi="Arial Rounded MT Bold"/><w:noProof/><w:lang w:val="es-PE"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251671552" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5C1D0D6A" wp14:editId="11086B1D"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:align>center</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>50800</wp:posOffset></wp:positionV><wp:extent cx="4343400" cy="1013460"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="550983624" name="Cuadro de texto 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="4343400" cy="1013460"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"><a:scrgbClr r="0" g="0" b="0"/></a:lnRef><a:fillRef idx="0"><a:scrgbClr r="0" g="0" b="0"/></a:fillRef><a:effectRef idx="0"><a:scrgbClr r="0" g="0" b="0"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="dk1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p w14:paraId="56B18094" w14:textId="40A5E851" w:rsidR="00CB7522" w:rsidRPr="00CB7522" w:rsidRDefault="00CB7522"><w:pPr><w:rPr><w:lang w:val="es-PE"/></w:rPr></w:pPr><w:r w:rsidRPr="000D1991"><w:rPr><w:lang w:val="es-PE"/></w:rPr><w:t>{{CRONOGRAMA}}</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shapetype w14:anchorId="5C1D0D6A" id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="Cuadro de texto 1" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;left:0;text-align:left;margin-left:0;margin-top:4pt;width:342pt;height:79.8pt;z-index:251671552;visibility:visible;mso-wrap-style:square;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:center;mso-position-horizontal-relative:page;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-height-percent:0;mso-height-relative:margin;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xA;IQAYJjozNwIAAH0EAAAOAAAAZHJzL2Uyb0RvYy54bWysVN+P2jAMfp+0/yHK+2gLHNsqyolxYpqE&#xA;7k7ipnsOaUKrpXGWBFr2188J5ddtT9OEFOzY+Wx/tju97xpF9sK6GnRBs0FKidAcylpvC/r9Zfnh&#xA;EyXOM10yBVoU9CAcvZ+9fzdtTS6GUIEqhSUIol3emoJW3ps8SRyvRMPcAIzQaJRgG+ZRtduktKxF&#xA;9EYlwzSdJC3Y0ljgwjm8fTga6SziSym4f5LSCU9UQTE3H08bz004k9mU5VvLTFXzPg32D1k0rNYY&#xA;9Az1wDwjO1v/AdXU3IID6QccmgSkrLmINWA1WfqmmnXFjIi1IDnOnGly/w+WP+7X5tkS332BDhsY&#xA;CGmNyx1ehno6aZvwj5kStCOFhzNtovOE4+V4hL8UTRxtWZqNxpNIbHJ5bqzzXwU0JAgFtdiXSBfb&#xA;r5zHkOh6cgnRHKi6XNZKRSXMglgoS/YMu6h8TBJf3HgpTdqCTkZ3aQS+sQXo8/uNYvxHKPMWATWl&#xA;8fJSfJB8t+l6RjZQHpAoC8cZcoYva8RdMeefmcWhQQJwEfwTHlIBJgO9REkF9tff7oM/9hKtlLQ4&#xA;hAV1P3fMCkrUN41d/pyNx2FqozK++zhExV5bNtcWvWsWgAxluHKGRzH4e3USpYXmFfdlHqKiiWmO&#xA;sQvqT+LCH1cD942L+Tw64Zwa5ld6bXiADh0JfL50r8yavp8eR+ERTuPK8jdtPfqGlxrmOw+yjj0P&#xA;BB9Z7XnHGY9t6fcxLNG1Hr0uX43ZbwAAAP//AwBQSwMEFAAGAAgAAAAhAJpcJsnZAAAABgEAAA8A&#xA;AABkcnMvZG93bnJldi54bWxMj0FPwzAMhe9I/IfISNxYCkIllKYToMGFExvi7DVZUtE4VZJ15d9j&#xA;TnCyrff0/L12vYRRzDblIZKG61UFwlIfzUBOw8fu5UqByAXJ4BjJavi2Gdbd+VmLjYknerfztjjB&#xA;IZQb1OBLmRopc+9twLyKkyXWDjEFLHwmJ03CE4eHUd5UVS0DDsQfPE722dv+a3sMGjZP7t71CpPf&#xA;KDMM8/J5eHOvWl9eLI8PIIpdyp8ZfvEZHTpm2scjmSxGDVykaFA8WKzVLS97dtV3Nciulf/xux8A&#xA;AAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250&#xA;ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAAACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAv&#xA;AQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAAACEAGCY6MzcCAAB9BAAADgAAAAAAAAAAAAAAAAAu&#xA;AgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYACAAAACEAmlwmydkAAAAGAQAADwAAAAAAAAAAAAAA&#xA;AACRBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAAAAAEAAQA8wAAAJcFAAAAAA==&#xA;" filled="f" stroked="f"><v:textbox><w:txbxContent><w:p w14:paraId="56B18094" w14:textId="40A5E851" w:rsidR="00CB7522" w:rsidRPr="00CB7522" w:rsidRDefault="00CB7522"><w:pPr><w:rPr><w:lang w:val="es-PE"/></w:rPr></w:pPr><w:r w:rsidRPr="000D1991"><w:rPr><w:lang w:val="es-PE"/></w:rPr><w:t>{{CRONOGRAMA}}</w:t></w:r></w:p></w:txbxContent></v:textbox><w10:wrap anchorx="page"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p>
'@

$hostRange.InsertXML($newParaXml)

Write-Host "Updated schedule text box: removed fill/outline, refreshed style refs"
